$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Weekly Progress" sheet: insert a new week's row at the top of the
#    data (row 2), pushing the existing 4 rows down by one.
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Progress")

$wsWeekly.Rows.Item(2).Insert()

# Copy the number format (date style) down from the row below so the new
# row's date cell keeps the same "m/d/yyyy" style used by the rest of the
# column, then fill in this week's values.
$wsWeekly.Cells.Item(3, 1).Copy()
$wsWeekly.Cells.Item(2, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsWeekly.Cells.Item(2, 1).Value = 43623
$wsWeekly.Cells.Item(2, 2).Value = 188
$wsWeekly.Cells.Item(2, 3).Formula = "=SUM(B2,-D2)"
$wsWeekly.Cells.Item(2, 4).Value = 86

# Grow the "Table15" structured table to include the new row.
$loWeekly = $wsWeekly.ListObjects.Item(1)
$loWeekly.Resize($wsWeekly.Range("A1:D6"))

# Update the chart's series so they reference the expanded ranges.
$chartObj = $wsWeekly.ChartObjects().Item(1)
$chart = $chartObj.Chart()

$series1 = $chart.SeriesCollection().Item(1)
$series1.Formula = "=SERIES('Weekly Progress'!`$C`$1,'Weekly Progress'!`$A`$2:`$A`$6,'Weekly Progress'!`$C`$2:`$C`$6,1)"

$series2 = $chart.SeriesCollection().Item(2)
$series2.Formula = "=SERIES('Weekly Progress'!`$D`$1,'Weekly Progress'!`$A`$2:`$A`$6,'Weekly Progress'!`$D`$2:`$D`$6,2)"

# The chart is anchored to cells; since a row was inserted above it, move
# it down by one row's worth of height so it keeps its visual position
# relative to the data beneath it.
$chartObj.Top = $chartObj.Top() + $wsWeekly.Rows.Item(2).RowHeight()

# Update the sheet's remembered selection.
$wsWeekly.Activate()
$wsWeekly.Range("C3").Select()

# ---------------------------------------------------------------------
# 2. "Areas Features Validations" sheet: mark three scenarios as "todo"
#    and tidy up a couple of row heights / the remembered selection.
# ---------------------------------------------------------------------
$wsAreas = $wb.Worksheets.Item("Areas Features Validations")

$wsAreas.Range("E164").Value = "todo"
$wsAreas.Range("E165").Value = "todo"
$wsAreas.Range("E166").Value = "todo"

$wsAreas.Rows.Item(87).AutoFit()
$wsAreas.Rows.Item(110).AutoFit()

$wsAreas.Activate()
$wsAreas.Range("E164:E166").Select()
